# Applies the "Add files via upload" revision to preguntas.xlsx:
#  - Removes the header row (item / pregunta / escala / posibles respuestas)
#    so the data now starts on row 1 and item numbers shift up by one.
#  - Converts the "escala" column (C) from descriptive text (likert/Binario)
#    to the numeric size of the scale (5, 2 or 3).
#  - Updates "posibles respuestas" (D) text for the items whose scale
#    changed to a 2-point or 3-point scale.
#  - Mirrors the left-alignment that a couple of the binary-scale rows
#    picked up in column D.
#  - Widens column B and re-selects the (now data) first row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old header row; rows 2..21 (items 1..20) shift up to 1..20.
$ws.Rows(1).Delete()

# 2) Column C ("escala"): replace the text label with the numeric count of
#    possible answers for that scale.
$scale = @(5, 2, 5, 2, 5, 2, 5, 5, 3, 5, 5, 5, 5, 5, 5, 5, 5, 5, 5, 3)
for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $scale[$i]
}

# 3) Column D ("posibles respuestas"): update the text for rows whose scale
#    is no longer the 5-point likert scale.
$binario = " 1: De acuerdo, 2: Totalmente de acuerdo"
$siNoNoSe = "1: SI, 2: NO, 3: NO SE"

$ws.Cells.Item(4, 4).Value = $binario
$ws.Cells.Item(6, 4).Value = $binario
$ws.Cells.Item(9, 4).Value = $siNoNoSe
$ws.Cells.Item(20, 4).Value = $siNoNoSe

# 4) Match the left alignment used on the binary-scale rows in column D.
$ws.Cells.Item(4, 4).HorizontalAlignment = -4131
$ws.Cells.Item(6, 4).HorizontalAlignment = -4131

# 5) Column B grew wider in the revision.
$ws.Columns(2).ColumnWidth = 93.25

# 6) Selection now highlights the (new) first row instead of D4.
$ws.Range("A1:XFD1").Select()

Write-Host "edit applied"
